# Appends newly-logged sensor readings to the three "…-LIFTER" sheets.
#
# Each sheet is an append-only log where every row shares one constant
# "packet template" (columns B..I) and only column A (an Excel date-time
# serial) changes row to row. Column A typically repeats the same serial a
# handful of times (duplicate packets captured back-to-back) before
# advancing to the next timestamp. The most-recently-appended row in the
# overall workbook is written before its timestamp has been converted to a
# proper date serial, so it still holds a plain "yyyy-mm-dd hh:mm:ss" text
# string with no date number-format applied; as soon as more rows land
# after it, that placeholder gets replaced with the real numeric serial
# (date-formatted) like all the others.
#
# This script re-plays that: it rewrites the tail of each sheet (promoting
# the previous "still text" row to a real numeric/date cell) and appends
# the freshly-captured rows, leaving the newest row of the fastest-moving
# sheet as the new "still text" placeholder.

function Set-SensorRow($ws, $row, $aValue, $template) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $aValue
    if (-not ($aValue -is [string])) {
        # Matches style index 2 in this workbook (numFmtId 165): the date
        # serial is displayed, not the raw float.
        $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
    $ws.Cells.Item($row, 2).Value = $template.B
    $ws.Cells.Item($row, 3).Value = $template.C
    $ws.Cells.Item($row, 4).Value = $template.D
    $ws.Cells.Item($row, 5).Value = $template.E
    $ws.Cells.Item($row, 6).Value = $template.F
    $ws.Cells.Item($row, 7).Value = $template.G
    $ws.Cells.Item($row, 8).Value = $template.H
    $ws.Cells.Item($row, 9).Value = $template.I
}

function Update-SensorSheet($wb, $sheetName, $rows, $template) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in ($rows.Keys | Sort-Object { [int]$_ })) {
        Set-SensorRow $ws $row $rows[$row] $template
    }
}

$wb = $excel.ActiveWorkbook

# ROW35-FE-LIFTER: rows 474-487 get the rows-484..487 values shifted back
# one "slot" (the trailing duplicate set grows from 2 to 3 repeats), row
# 488 graduates from placeholder text to a real date, and rows 489-496 are
# newly appended (finishing with a normal numeric row).
$sheet0Rows = @{
    474 = 45725.23111297454
    475 = 45725.23113506944
    476 = 45725.23113506944
    477 = 45725.23113506944
    478 = 45725.23113506944
    479 = 45725.23113506944
    480 = 45725.23115855324
    481 = 45725.23115855324
    482 = 45725.23115855324
    483 = 45725.23115855324
    484 = 45725.23115855324
    485 = 45725.73125641204
    486 = 45725.73125641204
    487 = 45725.73125641204
    488 = 45725.73127832176
    489 = 45725.73127832176
    490 = 45725.73127832176
    491 = 45725.73130123843
    492 = 45725.73130123843
    493 = 45725.73130123843
    494 = 45726.23139893518
    495 = 45726.23142038195
    496 = 45726.23144357639
}
$sheet0Template = @{
    B = "0x01,0x90"
    C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
    D = "0x01,0x90,"
    E = "0xd"
    F = 400
    G = 568631262647113769549824.0
    H = 400
    I = 13
}
Update-SensorSheet $wb "ROW35-FE-LIFTER" $sheet0Rows $sheet0Template

# ROW35-MID-LIFTER: same shape, one slot longer; its freshest row (511)
# hasn't been converted yet, so it stays literal placeholder text.
$sheet1Rows = @{
    488 = 45725.07909302083
    489 = 45725.07909302083
    490 = 45725.07909302083
    491 = 45725.07909302083
    492 = 45725.07911518519
    493 = 45725.07911518519
    494 = 45725.07913833333
    495 = 45725.07913833333
    496 = 45725.07913833333
    497 = 45725.07913833333
    498 = 45725.07913833333
    499 = 45725.57923533564
    500 = 45725.57923533564
    501 = 45725.57923533564
    502 = 45725.57925716435
    503 = 45725.57925716435
    504 = 45725.57925716435
    505 = 45725.57928042824
    506 = 45725.57928042824
    507 = 45725.57928042824
    508 = 45726.07937777778
    509 = 45726.07939922454
    510 = 45726.07942256945
    511 = "2025-03-10 13:54:22"
}
$sheet1Template = @{
    B = "0x01,0x90"
    C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
    D = "0x01,0x90,"
    E = "0xe"
    F = 400
    G = 568631262647113769549824.0
    H = 400
    I = 14
}
Update-SensorSheet $wb "ROW35-MID-LIFTER" $sheet1Rows $sheet1Template

# ROW02-MID-LIFTER: same shape as sheet 0 (its new last row is numeric).
$sheet3Rows = @{
    492 = 45725.22919952546
    493 = 45725.22919952546
    494 = 45725.22919952546
    495 = 45725.22919952546
    496 = 45725.22922125
    497 = 45725.22922125
    498 = 45725.22922125
    499 = 45725.22922125
    500 = 45725.22922125
    501 = 45725.22924497685
    502 = 45725.22924497685
    503 = 45725.7293421412
    504 = 45725.7293421412
    505 = 45725.7293421412
    506 = 45725.72936453704
    507 = 45725.72936453704
    508 = 45725.72936453704
    509 = 45725.7293875
    510 = 45725.7293875
    511 = 45725.7293875
    512 = 45726.2294844213
    513 = 45726.22950657408
    514 = 45726.22953034722
}
$sheet3Template = @{
    B = "0x01,0x90"
    C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
    D = "0x01,0x90,"
    E = "0x3"
    F = 400
    G = 568631262647113769549824.0
    H = 400
    I = 3
}
Update-SensorSheet $wb "ROW02-MID-LIFTER" $sheet3Rows $sheet3Template
